# The sheet had a stray/empty leading column A (the real data lived in
# columns B:F). Delete it so the data shifts left into A:E and the sheet
# dimension shrinks from A1:F3 to A1:E3. Deleting column A also moves the
# bordered/bold header style from B1:F1 to A1:E1, and drops the leftover
# border style that used to sit on A2:A3, matching the target layout.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns.Item(1).Delete()

# Fix the header label typo: "MODEL_CONDITION" -> "MODELCONDITION".
# After the column shift this header now lives in D1.
$ws.Range("D1").Value = "MODELCONDITION"
